# Transaction_List.xlsx — "Add download instructional search results"
#
# The committed header row (A1:S1) that labeled each column
# (Library Unit, Date of Consultation, Staff Pennkey, ...) is wiped so the
# sheet ships as a blank downloadable template. Clearing the contents
# (not deleting the cells) preserves the existing header styles/borders
# and row height while dropping the text - this also empties the shared
# string table since those headers were the only strings in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wipe the header labels in row 1 (A1:S1) but keep their formatting.
$ws.Range("A1:S1").ClearContents()

# Move the selection/view the way the refreshed template ships: scrolled
# over to column X with the cursor sitting on X1.
$ws.Range("X1").Select() | Out-Null
